$d = $word.ActiveDocument

# 1) Fix typo: "for saveing user score" -> "for saving user score"
#    (kept as its own run, sitting right after the " // " run, so use
#    InsertXML on the exact matched range instead of Find.Execute's
#    Replace, which would merge it with the neighboring identical-format
#    run.)
$typoXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>for saving user score</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$typoSearch = $d.Content.Duplicate
if ($typoSearch.Find.Execute("for saveing user score")) {
    $typoRange = $d.Range($typoSearch.Start, $typoSearch.End)
    $typoRange.InsertXML($typoXml) | Out-Null
}

# 2) Collapse "Display game " + "with START button" into one run with new wording
$d.Content.Find.Execute("Display game with START button", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Build answer buttons and add event listeners", 2) | Out-Null

# 3) Rebuild the tail of the runGame() section:
#    - drop the leading tab before "Remove START button from DOM"
#    - add new paragraphs for displaying the high score and looping over
#      the question array, wiring up the answer buttons
#    - replace the two trailing blank paragraphs with a tab-only paragraph
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("Remove START button from DOM")
if ($found) {
    $startPara = $searchRange.Paragraphs(1)
    $startPos = $startPara.Range.Start
    $endPos = $d.Paragraphs($d.Paragraphs.Count).Range.End
    $target = $d.Range($startPos, $endPos)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Remove START button from DOM</w:t></w:r></w:p><w:p><w:r><w:t>Display</w:t></w:r><w:r><w:t xml:space="preserve"> last high score and user initials from localStorage</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>for each question in qArr[]</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">call nextQuestion to assign fields to </w:t></w:r><w:r><w:t>buttons</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:t>call configureButton the style of each button and attach event listener</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $target.InsertXML($xml) | Out-Null
}
